$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168878197669983
$ws.Range("B1").Value = 2.439692974090576
$ws.Range("D1").Value = 2.366775751113892
$ws.Range("E1").Value = 1.233752489089966
